$wb = $excel.ActiveWorkbook

# Update the choice_filter formula on the "survey" sheet
$survey = $wb.Worksheets.Item("survey")
$survey.Range("I18").Value = "context.region === data('region')"

# Rename settings keys to use underscores on the "settings" sheet
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "form_id"
$settings.Range("A3").Value = "form_version"
$settings.Range("A4").Value = "form_title"
